$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 160, shifting existing rows 160-239 down to 161-240
$ws.Rows.Item(160).Insert()

# Populate the new row 160 by copying the (now shifted) original row 161,
# then overwriting the cells that actually differ in the new record.
$ws.Cells.Item(160, 1).Value = $ws.Cells.Item(161, 1).Value()
$ws.Cells.Item(160, 2).Value = $ws.Cells.Item(161, 2).Value()
$ws.Cells.Item(160, 3).Value = $ws.Cells.Item(161, 3).Value()
$ws.Cells.Item(160, 4).Value = 44636
$ws.Cells.Item(160, 5).Value = $ws.Cells.Item(161, 5).Value()
$ws.Cells.Item(160, 6).Value = $ws.Cells.Item(161, 6).Value()
$ws.Cells.Item(160, 7).Value = $ws.Cells.Item(161, 7).Value()
$ws.Cells.Item(160, 8).Value = $ws.Cells.Item(161, 8).Value()
$ws.Cells.Item(160, 9).Value = $ws.Cells.Item(161, 9).Value()
$ws.Cells.Item(160, 10).Value = 200
$ws.Cells.Item(160, 11).Value = 900
$ws.Cells.Item(160, 12).Value = 950
$ws.Cells.Item(160, 13).Value = 925
$ws.Cells.Item(160, 14).Value = $ws.Cells.Item(161, 14).Value()
$ws.Cells.Item(160, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(160, 16).Value = 925
$ws.Cells.Item(160, 17).Value = $ws.Cells.Item(161, 17).Value()
$ws.Cells.Item(160, 18).Value = $ws.Cells.Item(161, 18).Value()

# Match the date cell's number format/style (s="2" on column D) to the rest of column D
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(161, 4).NumberFormat()
